$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell I1 (AutoFill1)
$ws.Range("I1").Value = "AutoFill1"

# D2 should become the *text* string "1.1" rather than the numeric 1.1.
# A direct .Value assignment of "1.1" gets auto-coerced back to a number by
# Excel, so instead compute it as a text formula and then freeze it to a
# plain value via copy / paste-special (values only) so it ends up stored
# as a shared string, matching a manually-retyped text cell.
$ws.Range("D2").Formula = "=""1.1"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

# Corrected enum values that POI had previously written incorrectly.
$ws.Range("F7").Value = "Enum1Val2"
$ws.Range("F8").Value = "Enum1Val1"
$ws.Range("F9").Value = "Enum1Val1"

# G9's boolean formula result picks up the plain default style (no longer
# the separate, duplicate cellXf) - copy A1's (default) formatting onto it
# without touching its value/formula.
$ws.Range("A1").Copy()
$ws.Range("G9").PasteSpecial(-4122)

# New AutoFill1 output cell.
$ws.Range("E12").Value = "OK"

$excel.CutCopyMode = 0
$ws.Range("E12").Select()
